$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 13 (Docentes responsáveis: value "5840820 - Gustavo Aristides Santana
# Martinez" in B13/C13) is removed entirely; everything below shifts up one row.
$ws.Rows.Item(13).Delete()

# Row 10 (Objetivos:): value becomes the professor's name instead of the long
# "Ensinar a linguagem..." objectives text.
$ws.Range("B10").Value = "5840820 - Gustavo Aristides Santana Martinez"
$ws.Range("C10").Value = "5840820 - Gustavo Aristides Santana Martinez"

# Row 13 (now "Programa resumido:"): value becomes "Semestral" instead of the long
# summary-program text.
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 (now "Programa:"): value becomes the activation date instead of the long
# program text.
$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "01/01/2018"

# Row 18 (now "Método:"): value becomes the professor's name.
$ws.Range("B18").Value = "5840820 - Gustavo Aristides Santana Martinez"
$ws.Range("C18").Value = "5840820 - Gustavo Aristides Santana Martinez"

# Row 21 (now "Bibliografia:"): value becomes the "recuperação" norm text instead of
# the long bibliography text.
$ws.Range("B21").Value = "- A recuperação deverá consistir de uma prova englobando a matéria toda do semestre.- A média final (pós-recuperação) deverá ser composta por uma média simples entre a nota do semestre (nota final) e a da prova de recuperação."
$ws.Range("C21").Value = "- A recuperação deverá consistir de uma prova englobando a matéria toda do semestre.- A média final (pós-recuperação) deverá ser composta por uma média simples entre a nota do semestre (nota final) e a da prova de recuperação."
